$d = $word.ActiveDocument

# Anchor 1: end of the "Total vendido:133600" paragraph's text (just before its
# paragraph mark). Everything from here through the end of the
# "...Ejemplo7.java, según sea el caso)" paragraph (inclusive of its own
# paragraph mark) must be removed, per the commit's instructions cleanup.
$anchor1 = $d.Content
$anchor1.Find.Execute("Total vendido:133600") | Out-Null
$startDel = $anchor1.End

# Anchor 2: the last paragraph to be removed ends with this text.
$anchor2 = $d.Content
$anchor2.Find.Execute("Ejemplo7.java, según sea el caso)") | Out-Null
$endDel = $anchor2.End + 1   # +1 to include that paragraph's own paragraph mark

# Delete the whole block. The emulated Range.Delete() only reliably removes a
# single character at a time when the span crosses paragraph marks, so walk
# it one character at a time (the start offset stays fixed: every delete
# shifts the remaining text left by one character).
$total = $endDel - $startDel
for ($i = 0; $i -lt $total; $i++) {
    $chunk = $d.Range($startDel, $startDel + 1)
    $chunk.Delete()
}
